$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "StoreId" column (E) mirroring the existing D (Date) column's values.
$ws.Range("E1").Value = "StoreId"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 5
$ws.Range("E7").Value = 6
$ws.Range("E8").Value = 7
$ws.Range("E9").Value = 8
$ws.Range("E10").Value = 9
$ws.Range("E11").Value = 10

# Move the active selection the way the author's Excel session left it.
$ws.Range("E12").Select() | Out-Null
